$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily records to append (row, date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.)
$newRows = @(
    @(465, 44539, 1, 14, 232.9063383796373),
    @(466, 44540, 1, 12, 199.6340043254034),
    @(467, 44541, 3, 15, 249.5425054067543),
    @(468, 44542, 6, 17, 282.8148394609882),
    @(469, 44543, 12, 26, 432.5403427050408),
    @(470, 44544, 3, 26, 432.5403427050408),
    @(471, 44545, 0, 26, 432.5403427050408),
    @(472, 44546, 2, 27, 449.1765097321577),
    @(473, 44547, 3, 29, 482.4488437863916),
    @(474, 44548, 11, 37, 615.5381800033272),
    @(475, 44550, 7, 38, 632.1743470304442),
    @(476, 44551, 3, 29, 482.4488437863916),
    @(477, 44552, 0, 26, 432.5403427050408),
    @(478, 44553, 2, 28, 465.8126767592747),
    @(479, 44554, 3, 29, 482.4488437863916),
    @(480, 44555, 3, 29, 482.4488437863916),
    @(481, 44556, 1, 19, 316.0871735152221),
    @(482, 44557, 0, 12, 199.6340043254034),
    @(483, 44558, 8, 17, 282.8148394609882),
    @(484, 44559, 2, 19, 316.0871735152221),
    @(485, 44560, 6, 23, 382.6318416236899),
    @(486, 44561, 8, 28, 465.8126767592747),
    @(487, 44562, 7, 32, 532.3573448677425),
    @(488, 44563, 12, 43, 715.355182166029),
    @(489, 44564, 14, 57, 948.2615205456664),
    @(490, 44565, 5, 54, 898.3530194643154),
    @(491, 44566, 14, 66, 1097.987023789719)
)

$lastExistingRow = 464
$firstNewRow = 465
$lastNewRow = 491

# Copy the formatting (styles, number formats, borders, etc.) of the last
# existing data row down onto the new rows before filling in values.
$srcFormat = $ws.Range($ws.Cells.Item($lastExistingRow, 1), $ws.Cells.Item($lastExistingRow, 4))
$dstFormat = $ws.Range($ws.Cells.Item($firstNewRow, 1), $ws.Cells.Item($lastNewRow, 4))
$srcFormat.Copy()
$dstFormat.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

foreach ($r in $newRows) {
    $rowIndex = $r[0]
    $ws.Cells.Item($rowIndex, 1).Value = $r[1]
    $ws.Cells.Item($rowIndex, 2).Value = $r[2]
    $ws.Cells.Item($rowIndex, 3).Value = $r[3]
    $ws.Cells.Item($rowIndex, 4).Value = $r[4]
}
